# Refatorando o consolidador para modelo ETL
# Update absenteeism data rows 2-11 with new source values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = 30304; B = "Stella Moraes";        C = "Vendas";     D = "Outros";              E = 4; F = 45095; G = 10293.32 }
    @{ Row = 3;  A = 25915; B = "Sabrina Jesus";         C = "P&D";        D = "Problemas pessoais";   E = 3; F = 45079; G = 10529.68 }
    @{ Row = 4;  A = 34306; B = "Natália Cavalcanti";    C = "P&D";        D = "Viagem de negócios";   E = 8; F = 45103; G = 5448.32 }
    @{ Row = 5;  A = 86775; B = "Leandro Cunha";         C = "P&D";        D = "Viagem de negócios";   E = 3; F = 45094; G = 3095.09 }
    @{ Row = 6;  A = 65918; B = "Alexandre Oliveira";    C = "Vendas";     D = "Outros";              E = 6; F = 45086; G = 4947.71 }
    @{ Row = 7;  A = 76078; B = "Maria Fernanda Gomes";  C = "P&D";        D = "Viagem de negócios";   E = 1; F = 45081; G = 12473.04 }
    @{ Row = 8;  A = 64818; B = "João Pedro Sales";      C = "Vendas";     D = "Outros";              E = 5; F = 45105; G = 3776.71 }
    @{ Row = 9;  A = 24047; B = "Luiz Miguel Alves";     C = "Vendas";     D = "Viagem de negócios";   E = 3; F = 45099; G = 7914.07 }
    @{ Row = 10; A = 28185; B = "Dr. Davi Melo";         C = "Operações";  D = "Viagem de negócios";   E = 5; F = 45096; G = 5581.62 }
    @{ Row = 11; A = 6416;  B = "Pietro Nunes";          C = "Marketing";  D = "Doença";               E = 2; F = 45091; G = 4741.84 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
